$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "Done"
$ws.Range("H4").Value = "Done"
$ws.Range("H5").Value = "Done"
$ws.Range("H10").Value = "Done"
$ws.Range("H11").Value = "Done"
$ws.Range("H13").Value = "Done"
$ws.Range("H14").Value = "Done"
$ws.Range("H19").Value = "Done"
$ws.Range("H20").Value = "Done"
$ws.Range("H22").Value = "Done"
$ws.Range("H24").Value = "Done"

$ws.Range("I15").Select()
